$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.287.07'
$ws.Range("E2").Value = '  +2.99%  '
$ws.Range("D3").Value = '2.525.76'
$ws.Range("E3").Value = '  +3.77%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '535.82'
$ws.Range("E5").Value = '  +4.36%  '
$ws.Range("D6").Value = '135.83'
$ws.Range("E6").Value = '  +4.39%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.52%  '
$ws.Range("D8").Value = '0.567'
$ws.Range("E8").Value = '  +3.26%  '
$ws.Range("D9").Value = '2.525.99'
$ws.Range("E9").Value = '  +3.21%  '
$ws.Range("D10").Value = '0.0998'
$ws.Range("E10").Value = '  +4.93%  '
$ws.Range("E11").Value = '  -1.25%  '
$ws.Range("D12").Value = '5.23'
$ws.Range("E12").Value = '  +0.94%  '
$ws.Range("D13").Value = '0.336'
$ws.Range("E13").Value = '  +0.76%  '
$ws.Range("D14").Value = '2.971.66'
$ws.Range("E14").Value = '  +3.63%  '
$ws.Range("D15").Value = '59.233.76'
$ws.Range("E15").Value = '  +3.10%  '
$ws.Range("D16").Value = '22.51'
$ws.Range("E16").Value = '  +3.07%  '
$ws.Range("D17").Value = '0.0000137'
$ws.Range("E17").Value = '  +3.87%  '
$ws.Range("D18").Value = '2.527.01'
$ws.Range("E18").Value = '  +3.60%  '
$ws.Range("D19").Value = '10.76'
$ws.Range("E19").Value = '  +2.98%  '
$ws.Range("D20").Value = '4.25'
$ws.Range("E20").Value = '  +3.58%  '
$ws.Range("D21").Value = '323.71'
$ws.Range("E21").Value = '  +2.60%  '
$ws.Range("D22").Value = '6.17'
$ws.Range("E22").Value = '  +8.68%  '
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.18%  '
$ws.Range("D24").Value = '65.97'
$ws.Range("E24").Value = '  +3.85%  '
$ws.Range("D25").Value = '0.410'
$ws.Range("E25").Value = '  +0.65%  '
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  +0.28%  '
$ws.Range("D27").Value = '0.161'
$ws.Range("E27").Value = '  +1.38%  '
$ws.Range("D28").Value = '7.55'
$ws.Range("E28").Value = '  +4.97%  '
$ws.Range("D29").Value = '0.0₃0768'
$ws.Range("E29").Value = '  +6.32%  '
$ws.Range("D30").Value = '173.48'
$ws.Range("E30").Value = '  +1.69%  '
$ws.Range("E31").Value = '  +5.46%  '
$ws.Range("D32").Value = '1.22'
$ws.Range("E32").Value = '  +4.74%  '
$ws.Range("D33").Value = '6.39'
$ws.Range("E33").Value = '  +2.42%  '
$ws.Range("E34").Value = '  +0.07%  '
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.30%  '
$ws.Range("D36").Value = '18.26'
$ws.Range("E36").Value = '  +3.05%  '
$ws.Range("D37").Value = '1.26'
$ws.Range("E37").Value = '  -1.78%  '
$ws.Range("D38").Value = '4.01'
$ws.Range("E38").Value = '  +2.82%  '
$ws.Range("E39").Value = '  +4.86%  '
$ws.Range("D40").Value = '36.80'
$ws.Range("E40").Value = '  +1.67%  '
$ws.Range("D41").Value = '0.790'
$ws.Range("E41").Value = '  +1.49%  '
$ws.Range("D42").Value = '281.97'
$ws.Range("E42").Value = '  +5.07%  '
$ws.Range("D43").Value = '3.51'
$ws.Range("E43").Value = '  +4.10%  '
$ws.Range("D44").Value = '5.12'
$ws.Range("E44").Value = '  +4.29%  '
$ws.Range("D45").Value = '132.63'
$ws.Range("E45").Value = '  +10.33%  '
$ws.Range("D46").Value = '0.600'
$ws.Range("E46").Value = '  +2.14%  '
$ws.Range("E47").Value = '  +2.61%  '
$ws.Range("E48").Value = '  +5.86%  '
$ws.Range("E49").Value = '  +5.04%  '
$ws.Range("D50").Value = '17.26'
$ws.Range("E50").Value = '  +4.30%  '
$ws.Range("D51").Value = '1.767.64'
$ws.Range("E51").Value = '  +3.63%  '
